# Update LR-pair data per Dr Hou advice (Natmi rerun): now 3 expressing cells per group,
# and target cluster expands to both FAPs and sCs for each sending cluster.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Gnai2", "Lhcgr", "FAPs", 3, 1, 106.8785896666667, 320.635769, 0.1508748302900445, 0.1508748302900445, 1, 0.3333333333333333, 0.2064483333333333, 0.619345, 0.3625320258231888, 0.3625320258231888, 22.06490670570055, 198.584160351305, 0.05469695787077965, 0.05469695787077965),
    @("ECs", "Gnai2", "Lhcgr", "sCs", 3, 1, 106.8785896666667, 320.635769, 0.1508748302900445, 0.1508748302900445, 3, 1, 0.363014, 1.089042, 0.6374679741768112, 0.6374679741768112, 38.79842434925533, 349.185819143298, 0.09617787241926488, 0.09617787241926488),
    @("FAPs", "Gnai2", "Lhcgr", "FAPs", 3, 1, 77.232923, 231.698769, 0.1090256166999485, 0.1090256166999485, 1, 0.3333333333333333, 0.2064483333333333, 0.619345, 0.3625320258231888, 0.3625320258231888, 15.94460823181167, 143.501474086305, 0.03952527768885482, 0.03952527768885482),
    @("FAPs", "Gnai2", "Lhcgr", "sCs", 3, 1, 77.232923, 231.698769, 0.1090256166999485, 0.1090256166999485, 3, 1, 0.363014, 1.089042, 0.6374679741768112, 0.6374679741768112, 28.036632309922, 252.329690789298, 0.06950033901109369, 0.0695003390110937),
    @("M1", "Gnai2", "Lhcgr", "FAPs", 3, 1, 136.676337, 410.029011, 0.1929387280825172, 0.1929387280825172, 1, 0.3333333333333333, 0.2064483333333333, 0.619345, 0.3625320258231888, 0.3625320258231888, 28.216601979755, 253.949417817795, 0.06994646795150435, 0.06994646795150435),
    @("M1", "Gnai2", "Lhcgr", "sCs", 3, 1, 136.676337, 410.029011, 0.1929387280825172, 0.1929387280825172, 3, 1, 0.363014, 1.089042, 0.6374679741768112, 0.6374679741768112, 49.61542379971801, 446.538814197462, 0.1229922601310129, 0.1229922601310129),
    @("M2", "Gnai2", "Lhcgr", "FAPs", 3, 1, 143.4557186666667, 430.367156, 0.2025088212285795, 0.2025088212285795, 1, 0.3333333333333333, 0.2064483333333333, 0.619345, 0.3625320258231888, 0.3625320258231888, 29.61619402586889, 266.54574623282, 0.07341593320706293, 0.07341593320706291),
    @("M2", "Gnai2", "Lhcgr", "sCs", 3, 1, 143.4557186666667, 430.367156, 0.2025088212285795, 0.2025088212285795, 3, 1, 0.363014, 1.089042, 0.6374679741768112, 0.6374679741768112, 52.07643425606133, 468.687908304552, 0.1290928880215166, 0.1290928880215166),
    @("Neutro", "Gnai2", "Lhcgr", "FAPs", 3, 1, 198.5982106666667, 595.794632, 0.2803505493821544, 0.2803505493821544, 1, 0.3333333333333333, 0.2064483333333333, 0.619345, 0.3625320258231888, 0.3625320258231888, 41.00026959511555, 369.00242635604, 0.1016360526081564, 0.1016360526081563),
    @("Neutro", "Gnai2", "Lhcgr", "sCs", 3, 1, 198.5982106666667, 595.794632, 0.2803505493821544, 0.2803505493821544, 3, 1, 0.363014, 1.089042, 0.6374679741768112, 0.6374679741768112, 72.09393084694933, 648.845377622544, 0.178714496773998, 0.178714496773998),
    @("sCs", "Gnai2", "Lhcgr", "FAPs", 3, 1, 45.55066433333334, 136.651993, 0.06430145431675577, 0.06430145431675577, 1, 0.3333333333333333, 0.2064483333333333, 0.619345, 0.3625320258231888, 0.3625320258231888, 9.40385873384278, 84.634728604585, 0.0233113364968307, 0.0233113364968307),
    @("sCs", "Gnai2", "Lhcgr", "sCs", 3, 1, 45.55066433333334, 136.651993, 0.06430145431675577, 0.06430145431675577, 3, 1, 0.363014, 1.089042, 0.6374679741768112, 0.6374679741768112, 16.53552886230067, 148.819759760706, 0.04099011781992507, 0.04099011781992507)
)

$r = 2
foreach ($row in $data) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}